$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "GetAll" contacts return example (row 9, column D) to reflect
# the new GetAllContacts() return shape (includes alias per contact).
$ws.Range("D9").Value = "[{Status=Success}, {Username=user_1, alias=abc},{Username=user_2, alias=blabla},…]"

# Update the view: move the active selection to D10 instead of D15
# (this also brings the top-left visible cell back to A1, removing the
# previous scrolled-down A6 anchor).
$ws.Range("D10").Select()
